$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  "2" = @{
    "E" = 3
    "G" = 173.5887273333334
    "H" = 520.7661820000001
    "I" = 0.2624583749605043
    "J" = 0.2624583749605043
    "K" = 3
    "M" = 3.155977333333333
    "N" = 9.467932
    "O" = 0.3579027849973545
    "P" = 0.3579027849973545
    "Q" = 547.8420887861805
    "R" = 4930.578799075624
    "S" = 0.09393458334424443
    "T" = 0.0939345833442444
  }
  "3" = @{
    "E" = 3
    "G" = 173.5887273333334
    "H" = 520.7661820000001
    "I" = 0.2624583749605043
    "J" = 0.2624583749605043
    "K" = 3
    "M" = 3.165953666666667
    "N" = 9.497861
    "O" = 0.359034148472735
    "P" = 0.359034148472735
    "Q" = 549.573867792967
    "R" = 4946.164810136703
    "S" = 0.09423151916348246
    "T" = 0.09423151916348244
  }
  "4" = @{
    "E" = 3
    "G" = 173.5887273333334
    "H" = 520.7661820000001
    "I" = 0.2624583749605043
    "J" = 0.2624583749605043
    "K" = 3
    "M" = 2.496042666666666
    "N" = 7.488128
    "O" = 0.2830630665299106
    "P" = 0.2830630665299106
    "Q" = 433.2848698763663
    "R" = 3899.563828887296
    "S" = 0.07429227245277747
    "T" = 0.07429227245277746
  }
  "5" = @{
    "E" = 3
    "G" = 141.053299
    "H" = 423.159897
    "I" = 0.2132662656560029
    "J" = 0.2132662656560029
    "K" = 3
    "M" = 3.155977333333333
    "N" = 9.467932
    "O" = 0.3579027849973545
    "P" = 0.3579027849973545
    "Q" = 445.1610144358893
    "R" = 4006.449129923004
    "S" = 0.0763285904242691
    "T" = 0.07632859042426908
  }
  "6" = @{
    "E" = 3
    "G" = 141.053299
    "H" = 423.159897
    "I" = 0.2132662656560029
    "J" = 0.2132662656560029
    "K" = 3
    "M" = 3.165953666666667
    "N" = 9.497861
    "O" = 0.359034148472735
    "P" = 0.359034148472735
    "Q" = 446.5682091644797
    "R" = 4019.113882480317
    "S" = 0.07656987208776309
    "T" = 0.07656987208776309
  }
  "7" = @{
    "E" = 3
    "G" = 141.053299
    "H" = 423.159897
    "I" = 0.2132662656560029
    "J" = 0.2132662656560029
    "K" = 3
    "M" = 2.496042666666666
    "N" = 7.488128
    "O" = 0.2830630665299106
    "P" = 0.2830630665299106
    "Q" = 352.0750525780907
    "R" = 3168.675473202816
    "S" = 0.06036780314397076
    "T" = 0.06036780314397075
  }
  "8" = @{
    "E" = 3
    "G" = 322.1880443333333
    "H" = 966.564133
    "I" = 0.4871338815973437
    "J" = 0.4871338815973436
    "K" = 3
    "M" = 3.155977333333333
    "N" = 9.467932
    "O" = 0.3579027849973545
    "P" = 0.3579027849973545
    "Q" = 1016.818164986995
    "R" = 9151.363484882955
    "S" = 0.1743465728902608
    "T" = 0.1743465728902608
  }
  "9" = @{
    "E" = 3
    "G" = 322.1880443333333
    "H" = 966.564133
    "I" = 0.4871338815973437
    "J" = 0.4871338815973436
    "K" = 3
    "M" = 3.165953666666667
    "N" = 9.497861
    "O" = 0.359034148472735
    "P" = 0.359034148472735
    "Q" = 1020.032420313279
    "R" = 9180.291782819513
    "S" = 0.1748976983715204
    "T" = 0.1748976983715204
  }
  "10" = @{
    "E" = 3
    "G" = 322.1880443333333
    "H" = 966.564133
    "I" = 0.4871338815973437
    "J" = 0.4871338815973436
    "K" = 3
    "M" = 2.496042666666666
    "N" = 7.488128
    "O" = 0.2830630665299106
    "P" = 0.2830630665299106
    "Q" = 804.1951053458914
    "R" = 7237.755948113024
    "S" = 0.1378896103355625
    "T" = 0.1378896103355625
  }
  "11" = @{
    "E" = 3
    "G" = 24.56519766666666
    "H" = 73.69559299999999
    "I" = 0.03714147778614916
    "J" = 0.03714147778614916
    "K" = 3
    "M" = 3.155977333333333
    "N" = 9.467932
    "O" = 0.3579027849973545
    "P" = 0.3579027849973545
    "Q" = 77.52720702485287
    "R" = 697.7448632236758
    "S" = 0.01329303833858016
    "T" = 0.01329303833858016
  }
  "12" = @{
    "E" = 3
    "G" = 24.56519766666666
    "H" = 73.69559299999999
    "I" = 0.03714147778614916
    "J" = 0.03714147778614916
    "K" = 3
    "M" = 3.165953666666667
    "N" = 9.497861
    "O" = 0.359034148472735
    "P" = 0.359034148472735
    "Q" = 77.77227762517477
    "R" = 699.9504986265729
    "S" = 0.01333505884996906
    "T" = 0.01333505884996906
  }
  "13" = @{
    "E" = 3
    "G" = 24.56519766666666
    "H" = 73.69559299999999
    "I" = 0.03714147778614916
    "J" = 0.03714147778614916
    "K" = 3
    "M" = 2.496042666666666
    "N" = 7.488128
    "O" = 0.2830630665299106
    "P" = 0.2830630665299106
    "Q" = 61.31578149110043
    "R" = 551.8420334199039
    "S" = 0.01051338059759994
    "T" = 0.01051338059759994
  }
}

foreach ($r in $data.Keys) {
  foreach ($c in $data[$r].Keys) {
    $ws.Range("$c$r").Value = $data[$r][$c]
  }
}